# Apply updated Betfair Back/Lay odds for 2026-01-08 (rows 2-12)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.38
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.94
$ws.Range("K2").Value = 3.15
$ws.Range("P2").Value = 1.43
$ws.Range("Q2").Value = 2.98
$ws.Range("S2").Value = 6.6
$ws.Range("T2").Value = 2.34
$ws.Range("V2").Value = 1.3
$ws.Range("W2").Value = 1.73
$ws.Range("X2").Value = 7.6
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 34
$ws.Range("AB2").Value = 6.6
$ws.Range("AC2").Value = 7.6
$ws.Range("AD2").Value = 24
$ws.Range("AE2").Value = 340
$ws.Range("AF2").Value = 15
$ws.Range("AG2").Value = 14.5
$ws.Range("AH2").Value = 36
$ws.Range("AJ2").Value = 38
$ws.Range("AK2").Value = 46

# Row 3
$ws.Range("F3").Value = 2.32
$ws.Range("G3").Value = 2.42
$ws.Range("I3").Value = 3.65
$ws.Range("J3").Value = 3.3
$ws.Range("K3").Value = 3.45
$ws.Range("L3").Value = 1.66
$ws.Range("M3").Value = 1.14
$ws.Range("O3").Value = 1.7
$ws.Range("P3").Value = 1.44
$ws.Range("Q3").Value = 2.98
$ws.Range("S3").Value = 6.8
$ws.Range("T3").Value = 2.34
$ws.Range("U3").Value = 1.65
$ws.Range("V3").Value = 1.38
$ws.Range("W3").Value = 1.7
$ws.Range("X3").Value = 7.6
$ws.Range("Y3").Value = 9.4
$ws.Range("Z3").Value = 23
$ws.Range("AA3").Value = 85
$ws.Range("AB3").Value = 6.6
$ws.Range("AC3").Value = 7.8
$ws.Range("AF3").Value = 13.5
$ws.Range("AG3").Value = 13
$ws.Range("AI3").Value = 120
$ws.Range("AJ3").Value = 38
$ws.Range("AK3").Value = 40
$ws.Range("AL3").Value = 200
$ws.Range("AM3").Value = 790
$ws.Range("AN3").Value = 600
$ws.Range("AO3").Value = 110

# Row 4
$ws.Range("G4").Value = 1.16
$ws.Range("H4").Value = 22
$ws.Range("I4").Value = 28
$ws.Range("L4").Value = 1.19
$ws.Range("N4").Value = 11
$ws.Range("P4").Value = 4.1
$ws.Range("R4").Value = 2.26
$ws.Range("T4").Value = 1.92
$ws.Range("U4").Value = 1.97
$ws.Range("W4").Value = 7.2
$ws.Range("Z4").Value = 290
$ws.Range("AB4").Value = 28
$ws.Range("AC4").Value = 30
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 970
$ws.Range("AK4").Value = 13
$ws.Range("AL4").Value = 48
$ws.Range("AN4").Value = 2.6

# Row 5
$ws.Range("F5").Value = 1.93
$ws.Range("G5").Value = 2.08
$ws.Range("J5").Value = 2.94
$ws.Range("N5").Value = 2.28
$ws.Range("Q5").Value = 3.1
$ws.Range("R5").Value = 1.14
$ws.Range("W5").Value = 1.94

# Row 6
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 3.4
$ws.Range("L6").Value = 1.6
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 2.78
$ws.Range("O6").Value = 1.54
$ws.Range("P6").Value = 1.58
$ws.Range("Q6").Value = 2.66
$ws.Range("S6").Value = 5.5
$ws.Range("T6").Value = 2.12
$ws.Range("U6").Value = 1.86
$ws.Range("V6").Value = 1.41
$ws.Range("X6").Value = 8.6
$ws.Range("Z6").Value = 20
$ws.Range("AB6").Value = 7.6
$ws.Range("AE6").Value = 50
$ws.Range("AK6").Value = 36
$ws.Range("AL6").Value = 65
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 40

# Row 7
$ws.Range("F7").Value = 1.6
$ws.Range("H7").Value = 5.3
$ws.Range("L7").Value = 1.26
$ws.Range("N7").Value = 6.8
$ws.Range("O7").Value = 1.15
$ws.Range("P7").Value = 2.92
$ws.Range("Q7").Value = 1.49
$ws.Range("R7").Value = 1.75
$ws.Range("S7").Value = 2.26
$ws.Range("T7").Value = 1.58
$ws.Range("U7").Value = 2.52
$ws.Range("X7").Value = 38
$ws.Range("Y7").Value = 38
$ws.Range("Z7").Value = 240
$ws.Range("AA7").Value = 140
$ws.Range("AB7").Value = 14.5
$ws.Range("AC7").Value = 12
$ws.Range("AD7").Value = 21
$ws.Range("AF7").Value = 13
$ws.Range("AG7").Value = 10.5
$ws.Range("AH7").Value = 17.5
$ws.Range("AI7").Value = 55
$ws.Range("AJ7").Value = 16.5
$ws.Range("AL7").Value = 25
$ws.Range("AM7").Value = 65
$ws.Range("AN7").Value = 5.8
$ws.Range("AO7").Value = 46

# Row 8
$ws.Range("F8").Value = 4.9
$ws.Range("H8").Value = 1.84
$ws.Range("I8").Value = 1.85
$ws.Range("J8").Value = 3.75
$ws.Range("K8").Value = 3.95
$ws.Range("L8").Value = 1.43
$ws.Range("N8").Value = 3.8
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 1.94
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.35
$ws.Range("S8").Value = 3.55
$ws.Range("T8").Value = 1.87
$ws.Range("V8").Value = 2.16
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 8.8
$ws.Range("Z8").Value = 11
$ws.Range("AA8").Value = 20
$ws.Range("AB8").Value = 17.5
$ws.Range("AC8").Value = 8.8
$ws.Range("AE8").Value = 20
$ws.Range("AF8").Value = 38
$ws.Range("AG8").Value = 19.5
$ws.Range("AH8").Value = 20
$ws.Range("AI8").Value = 40
$ws.Range("AJ8").Value = 120
$ws.Range("AK8").Value = 70
$ws.Range("AL8").Value = 75
$ws.Range("AM8").Value = 120
$ws.Range("AN8").Value = 75
$ws.Range("AO8").Value = 13

# Row 9
$ws.Range("F9").Value = 1.69
$ws.Range("G9").Value = 1.76
$ws.Range("H9").Value = 6.6
$ws.Range("I9").Value = 8
$ws.Range("K9").Value = 3.75
$ws.Range("N9").Value = 2.66
$ws.Range("O9").Value = 1.52
$ws.Range("P9").Value = 1.55
$ws.Range("Q9").Value = 2.58
$ws.Range("S9").Value = 5.1
$ws.Range("U9").Value = 1.61
$ws.Range("W9").Value = 2.3
$ws.Range("Y9").Value = 18
$ws.Range("AD9").Value = 32
$ws.Range("AF9").Value = 8.8
$ws.Range("AJ9").Value = 18.5
$ws.Range("AK9").Value = 25
$ws.Range("AN9").Value = 18.5

# Row 10
$ws.Range("I10").Value = 10
$ws.Range("L10").Value = 1.48
$ws.Range("N10").Value = 3.35
$ws.Range("P10").Value = 1.81
$ws.Range("Q10").Value = 2.2
$ws.Range("S10").Value = 4.2
$ws.Range("T10").Value = 2.42
$ws.Range("V10").Value = 1.11
$ws.Range("Y10").Value = 24
$ws.Range("AA10").Value = 420
$ws.Range("AF10").Value = 7.4
$ws.Range("AG10").Value = 10
$ws.Range("AI10").Value = 190
$ws.Range("AN10").Value = 9.4
$ws.Range("AO10").Value = 350

# Row 11
$ws.Range("F11").Value = 1.61
$ws.Range("G11").Value = 1.62
$ws.Range("H11").Value = 6.2
$ws.Range("N11").Value = 4.2
$ws.Range("O11").Value = 1.3
$ws.Range("P11").Value = 2.08
$ws.Range("Q11").Value = 1.9
$ws.Range("U11").Value = 2.02
$ws.Range("W11").Value = 2.6
$ws.Range("AC11").Value = 9.4
$ws.Range("AD11").Value = 24
$ws.Range("AH11").Value = 22
$ws.Range("AI11").Value = 80
$ws.Range("AN11").Value = 9

# Row 12
$ws.Range("F12").Value = 2.12
$ws.Range("G12").Value = 2.24
$ws.Range("I12").Value = 4.5
$ws.Range("J12").Value = 3.15
$ws.Range("K12").Value = 3.35
$ws.Range("N12").Value = 3.1
$ws.Range("O12").Value = 1.43
$ws.Range("P12").Value = 1.7
$ws.Range("Q12").Value = 2.34
$ws.Range("R12").Value = 1.26
$ws.Range("S12").Value = 4.4
$ws.Range("T12").Value = 1.9
$ws.Range("U12").Value = 1.96
$ws.Range("W12").Value = 1.8
$ws.Range("X12").Value = 11.5
$ws.Range("Y12").Value = 13
$ws.Range("Z12").Value = 30
$ws.Range("AA12").Value = 100
$ws.Range("AB12").Value = 8.2
$ws.Range("AC12").Value = 7.2
$ws.Range("AD12").Value = 17
$ws.Range("AE12").Value = 60
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 11
$ws.Range("AH12").Value = 20
$ws.Range("AI12").Value = 80
$ws.Range("AK12").Value = 26
$ws.Range("AL12").Value = 48
$ws.Range("AM12").Value = 170
$ws.Range("AN12").Value = 25
$ws.Range("AO12").Value = 80
